{"js": "// Replace two-digit division exercise answers in the table cells.\n// Each entry is unique plain text within a table-cell run, so a direct\n// body.search(...) + insertText(..., Word.InsertLocation.replace) keeps\n// the existing run formatting (font / size) intact.\nconst replacements = [\n  [\"59\u00f78=7, 3\", \"27\u00f75=5, 2\"],\n  [\"45\u00f77=6, 3\", \"42\u00f73=14, 0\"],\n  [\"34\u00f78=4, 2\", \"30\u00f74=7, 2\"],\n  [\"38\u00f72=19, 0\", \"38\u00f75=7, 3\"],\n  [\"50\u00f73=16, 2\", \"88\u00f75=17, 3\"],\n  [\"45\u00f75=9, 0\", \"23\u00f73=7, 2\"],\n  [\"47\u00f72=23, 1\", \"38\u00f79=4, 2\"],\n  [\"56\u00f73=18, 2\", \"37\u00f78=4, 5\"],\n  [\"10\u00f77=1, 3\", \"96\u00f79=10, 6\"],\n  [\"32\u00f76=5, 2\", \"46\u00f78=5, 6\"],\n  [\"26\u00f76=4, 2\", \"27\u00f72=13, 1\"],\n  [\"29\u00f75=5, 4\", \"24\u00f76=4, 0\"],\n  [\"44\u00f79=4, 8\", \"45\u00f74=11, 1\"],\n  [\"66\u00f72=33, 0\", \"37\u00f75=7, 2\"],\n  [\"14\u00f72=7, 0\", \"86\u00f72=43, 0\"],\n  [\"35\u00f76=5, 5\", \"18\u00f72=9, 0\"],\n  [\"79\u00f77=11, 2\", \"16\u00f77=2, 2\"],\n  [\"95\u00f74=23, 3\", \"83\u00f76=13, 5\"],\n  [\"77\u00f77=11, 0\", \"47\u00f74=11, 3\"],\n  [\"38\u00f77=5, 3\", \"12\u00f73=4, 0\"],\n  [\"41\u00f78=5, 1\", \"20\u00f74=5, 0\"],\n  [\"90\u00f78=11, 2\", \"82\u00f76=13, 4\"],\n  [\"90\u00f74=22, 2\", \"38\u00f79=4, 2\"],\n  [\"43\u00f77=6, 1\", \"91\u00f78=11, 3\"],\n  [\"23\u00f72=11, 1\", \"26\u00f77=3, 5\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace two-digit division exercise answers in the table cells.\n# Each old value is unique plain text within a single table-cell run, so\n# Find/Execute with a literal Replacement.Text keeps the run formatting\n# (font / size) untouched while only swapping the visible characters.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"59\u00f78=7, 3\", \"27\u00f75=5, 2\"),\n  @(\"45\u00f77=6, 3\", \"42\u00f73=14, 0\"),\n  @(\"34\u00f78=4, 2\", \"30\u00f74=7, 2\"),\n  @(\"38\u00f72=19, 0\", \"38\u00f75=7, 3\"),\n  @(\"50\u00f73=16, 2\", \"88\u00f75=17, 3\"),\n  @(\"45\u00f75=9, 0\", \"23\u00f73=7, 2\"),\n  @(\"47\u00f72=23, 1\", \"38\u00f79=4, 2\"),\n  @(\"56\u00f73=18, 2\", \"37\u00f78=4, 5\"),\n  @(\"10\u00f77=1, 3\", \"96\u00f79=10, 6\"),\n  @(\"32\u00f76=5, 2\", \"46\u00f78=5, 6\"),\n  @(\"26\u00f76=4, 2\", \"27\u00f72=13, 1\"),\n  @(\"29\u00f75=5, 4\", \"24\u00f76=4, 0\"),\n  @(\"44\u00f79=4, 8\", \"45\u00f74=11, 1\"),\n  @(\"66\u00f72=33, 0\", \"37\u00f75=7, 2\"),\n  @(\"14\u00f72=7, 0\", \"86\u00f72=43, 0\"),\n  @(\"35\u00f76=5, 5\", \"18\u00f72=9, 0\"),\n  @(\"79\u00f77=11, 2\", \"16\u00f77=2, 2\"),\n  @(\"95\u00f74=23, 3\", \"83\u00f76=13, 5\"),\n  @(\"77\u00f77=11, 0\", \"47\u00f74=11, 3\"),\n  @(\"38\u00f77=5, 3\", \"12\u00f73=4, 0\"),\n  @(\"41\u00f78=5, 1\", \"20\u00f74=5, 0\"),\n  @(\"90\u00f78=11, 2\", \"82\u00f76=13, 4\"),\n  @(\"90\u00f74=22, 2\", \"38\u00f79=4, 2\"),\n  @(\"43\u00f77=6, 1\", \"91\u00f78=11, 3\"),\n  @(\"23\u00f72=11, 1\", \"26\u00f77=3, 5\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  $rng.Find.Execute(\n    $oldText,   # FindText\n    $false,     # MatchCase\n    $false,     # MatchWholeWord\n    $false,     # MatchWildcards\n    $false,     # MatchSoundsLike\n    $false,     # MatchAllWordForms\n    $true,      # Forward\n    1,          # Wrap = wdFindContinue\n    $false,     # Format\n    $newText,   # ReplaceWith\n    2           # Replace = wdReplaceAll\n  )\n}\n"}
